# Apply changes described by the diff:
#  - Insert two new worksheets "TagServer" and "TagVolume" between "Tags" and
#    "AddRouteRules".
#  - Populate them with the server / volume tagging sample data.
#  - Format the header rows (bold font on a light-grey filled, thin-bordered
#    background) matching the look of the existing "Tags" sheet header.
#  - Make "TagVolume" the active / selected sheet and scroll the tab strip so
#    it (and its neighbours) are visible.

$wb = $excel.ActiveWorkbook

$tagsSheet = $wb.Worksheets.Item("Tags")

# --- Create "TagServer" right after "Tags" -------------------------------
$tagServer = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tagsSheet)
$tagServer.Name = "TagServer"

$tagServer.Range("A1").Value = "Hostname"
$tagServer.Range("B1").Value = "Application"
$tagServer.Range("C1").Value = "OS"
$tagServer.Range("D1").Value = "Subnet"
$tagServer.Range("E1").Value = "Enviroment"

$tagServer.Range("A2").Value = "Shruthi"
$tagServer.Range("B2").Value = "Demo=True"
$tagServer.Range("C2").Value = "Windows2012R2=True"
$tagServer.Range("D2").Value = "HubPublicSub=True"
$tagServer.Range("E2").Value = "dev=True"

$tagServer.Range("A3").Value = "Lakshmi"
$tagServer.Range("B3").Value = "Demo=True"
$tagServer.Range("C3").Value = "Windows2012R2=True"
$tagServer.Range("D3").Value = "HubPublicSub=True"
$tagServer.Range("E3").Value = "dev=True"

$tagsSheet.Range("B1").Copy() | Out-Null
$tagServer.Range("A1:E1").PasteSpecial(-4122) | Out-Null

# --- Create "TagVolume" right after "TagServer" --------------------------
$tagVolume = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tagServer)
$tagVolume.Name = "TagVolume"

$tagVolume.Range("A1").Value = "VolumeName"
$tagVolume.Range("B1").Value = "Application"
$tagVolume.Range("C1").Value = "OS"
$tagVolume.Range("D1").Value = "Subnet"
$tagVolume.Range("E1").Value = "Enviroment"

$tagVolume.Range("A2").Value = "Shruthi_disk2"
$tagVolume.Range("B2").Value = "Demo=True"
$tagVolume.Range("C2").Value = "Windows2012R2=True"
$tagVolume.Range("D2").Value = "HubPublicSub=True"
$tagVolume.Range("E2").Value = "dev=True"

$tagVolume.Range("A3").Value = "Lakshmi_disk2"
$tagVolume.Range("B3").Value = "Demo=True"
$tagVolume.Range("C3").Value = "Windows2012R2=True"
$tagVolume.Range("D3").Value = "HubPublicSub=True"
$tagVolume.Range("E3").Value = "dev=True"

$tagsSheet.Range("B1").Copy() | Out-Null
$tagVolume.Range("A1:E1").PasteSpecial(-4122) | Out-Null

$tagVolume.Columns("A:A").ColumnWidth = 13.85546875
$tagVolume.Columns("B:B").ColumnWidth = 11.28515625
$tagVolume.Columns("E:E").ColumnWidth = 11.28515625

# --- Selection / active-sheet bookkeeping --------------------------------
$tagServer.Range("A1:E3").Select() | Out-Null
$tagVolume.Range("A1:E3").Select() | Out-Null

$tagVolume.Activate()

Write-Host "done"
